$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "http://localhost:8080/api/v1/area/suburb/" + "{postcode}"
#    runs into a single run with the combined text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "http://localhost:8080/api/v1/area/suburb/{postcode}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "http://localhost:8080/api/v1/area/suburb/{postcode}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Mark the four "naked" inline image runs as <w:noProof/> (the ones that
#    did not already carry it).
# ---------------------------------------------------------------------------
$d.InlineShapes.Item(1).Range.NoProofing = $true
$d.InlineShapes.Item(2).Range.NoProofing = $true
$d.InlineShapes.Item(3).Range.NoProofing = $true
$d.InlineShapes.Item(4).Range.NoProofing = $true

# ---------------------------------------------------------------------------
# 3) Rewrite the "GET request needs to be sent twice ..." sentence, splitting
#    it into multiple runs and tweaking the wording.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "GET request needs to be sent twice to authenticate the user before using the cookie and CSRF token from the response to make a POST request.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$start = $rng.Start
$rng.Delete()

$pieces = @(
    "GET request needs to be sent ",
    "TWICE",
    " to authenticate the user before ",
    "being able to ",
    "us",
    "e",
    " the cookie and CSRF token from the response to make a POST request."
)

$cur = $start
foreach ($piece in $pieces) {
    $ins = $d.Range($cur, $cur)
    $ins.InsertAfter($piece)
    $cur = $cur + $piece.Length
}

# ---------------------------------------------------------------------------
# 4) Append the new "Additionally Note that JSESSIONID= ..." sentence after
#    the "X-CSRF-TOKEN." paragraph.
# ---------------------------------------------------------------------------
$tail = $d.Content
$tail.Find.Execute("X-CSRF-TOKEN.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertAt = $tail.End

$extraPieces = @(
    @{ text = " Additionally Note that "; highlight = $false },
    @{ text = "JSESSIONID="; highlight = $true },
    @{ text = " needs to be appended before pasting the cookie value. Please refer to the below image for adding the headers to the Request. "; highlight = $false }
)

$cur2 = $insertAt
foreach ($p in $extraPieces) {
    $ins = $d.Range($cur2, $cur2)
    $ins.InsertAfter($p.text)
    if ($p.highlight) {
        $hr = $d.Range($cur2, $cur2 + $p.text.Length)
        $hr.HighlightColorIndex = 7
    }
    $cur2 = $cur2 + $p.text.Length
}
